$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 30
$dstRow = 31

# Insert a new row below the last data row by copying the whole row above it,
# so the new row inherits exactly the same per-column cell styles used
# throughout the rest of the table.
$ws.Rows($srcRow).Copy()
$ws.Rows($dstRow).Insert()
$excel.CutCopyMode = 0

# Fill in the new record: 23/05/2018, Licata Rosa, Silesia Nera, Mt., 9
$ws.Cells.Item($dstRow, 1).Value = 43243
$ws.Cells.Item($dstRow, 2).Value = "Licata Rosa"
$ws.Cells.Item($dstRow, 3).Value = "Silesia Nera"
$ws.Cells.Item($dstRow, 4).Value = "Mt."
$ws.Cells.Item($dstRow, 5).Value = 9
